# Weekly update: insert a new price record as row 10, pushing the
# existing rows 10-19 down to rows 11-20 (new dimension A1:R20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 10, shifting rows 10:19
# (and any formatting) down to rows 11:20.
$ws.Rows.Item(10).EntireRow.Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44771
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112043
$ws.Range("G10").Value = "Pepino dulce"
$ws.Range("H10").Value = "Cultivar IV Región"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 17000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 17500
$ws.Range("N10").Value = "`$/bandeja 18 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 972
$ws.Range("Q10").Value = 18
$ws.Range("R10").Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by
# the rest of column D (style index 2 in the original workbook).
$ws.Range("D10").NumberFormat = $ws.Range("D11").NumberFormat
